$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # F4Disc
$ws2 = $wb.Worksheets.Item(2)   # F103 Robotex

# ---------------------------------------------------------------
# Sheet 1 (F4Disc): update prescale & BS1/BS2 inputs, add sample
# point % column (D14/D15)
# ---------------------------------------------------------------

# B6 prescale value: 41 -> 57
$ws1.Range("B6").Value = 57

# B7 formula now divides by 84 instead of 42
$ws1.Range("B7").Formula = "=(B6+1)*1/84"

# B10 (BS1 multiplier) 5 -> 7
$ws1.Range("B10").Value = 7

# B14 (BS2 multiplier) 2 -> 4
$ws1.Range("B14").Value = 4

# New label cell D14: "sample point %"
$ws1.Range("D14").Value = "sample point %"

# New formula cell D15: sample point ratio
$ws1.Range("D15").Formula = "=(B11+B7)/(B7+B11+B15)"

# ---------------------------------------------------------------
# Sheet 2 (F103 Robotex): add the same sample point % column
# (values/formulas on this sheet are otherwise unchanged)
# ---------------------------------------------------------------

$ws2.Range("D14").Value = "sample point %"
$ws2.Range("D15").Formula = "=(B11+B7)/(B7+B11+B15)"

# ---------------------------------------------------------------
# Update the active selections on each sheet to match the new
# editing focus, and leave sheet 1 as the selected tab (as before).
# ---------------------------------------------------------------

$ws2.Range("B26").Select()
$ws1.Range("B7").Select()
